# Scheduled runner refresh: update cached market-board / profit figures
# (currentAveragePrice*, Leve*Price*, Leve*Profit* columns) across all
# Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly
# pulled data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 30310906
$ws.Range("I43").Value = 83335210
$ws.Range("K43").Value = 83335210
$ws.Range("M43").Value = -83335141

$ws.Range("H86").Value = 71431816
$ws.Range("I86").Value = 100003020
$ws.Range("J86").Value = 3801
$ws.Range("K86").Value = 100003020
$ws.Range("L86").Value = 3801
$ws.Range("M86").Value = -100001897
$ws.Range("N86").Value = -6047

$ws.Range("H89").Value = 71431816
$ws.Range("I89").Value = 100003020
$ws.Range("J89").Value = 3801
$ws.Range("K89").Value = 500015100
$ws.Range("L89").Value = 19005
$ws.Range("M89").Value = -500009484
$ws.Range("N89").Value = -30237

$ws.Range("H97").Value = 2900
$ws.Range("J97").Value = 2900
$ws.Range("L97").Value = 8700
$ws.Range("N97").Value = -9692

$ws.Range("H112").Value = 126666.625
$ws.Range("J112").Value = 126666.625
$ws.Range("L112").Value = 379999.875
$ws.Range("N112").Value = -382215.875

$ws.Range("H118").Value = 839.6667
$ws.Range("I118").Value = 743
$ws.Range("J118").Value = 1033
$ws.Range("K118").Value = 2229
$ws.Range("L118").Value = 3099
$ws.Range("M118").Value = -572
$ws.Range("N118").Value = -6413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 55000
$ws.Range("J42").Value = 55000
$ws.Range("L42").Value = 55000
$ws.Range("N42").Value = -55972

$ws.Range("H61").Value = 1256056.9
$ws.Range("I61").Value = 4340.073
$ws.Range("J61").Value = 9809455
$ws.Range("K61").Value = 4340.073
$ws.Range("L61").Value = 9809455
$ws.Range("M61").Value = -4128.073
$ws.Range("N61").Value = -9809879

$ws.Range("H74").Value = 355576.2
$ws.Range("I74").Value = 1709.8125
$ws.Range("J74").Value = 720857.5600000001
$ws.Range("K74").Value = 1709.8125
$ws.Range("L74").Value = 720857.5600000001
$ws.Range("M74").Value = -835.8125
$ws.Range("N74").Value = -722605.5600000001

$ws.Range("H77").Value = 355576.2
$ws.Range("I77").Value = 1709.8125
$ws.Range("J77").Value = 720857.5600000001
$ws.Range("K77").Value = 8549.0625
$ws.Range("L77").Value = 3604287.8
$ws.Range("M77").Value = -4181.0625
$ws.Range("N77").Value = -3613023.8

$ws.Range("H92").Value = 59999
$ws.Range("J92").Value = 59999
$ws.Range("L92").Value = 59999
$ws.Range("N92").Value = -64991

$ws.Range("H97").Value = 6313.6113
$ws.Range("I97").Value = 6313.6113
$ws.Range("K97").Value = 6313.6113
$ws.Range("M97").Value = -5817.6113

$ws.Range("H102").Value = 2730.3635
$ws.Range("I102").Value = 2717.524
$ws.Range("K102").Value = 2717.524
$ws.Range("M102").Value = -1095.524

$ws.Range("H122").Value = 1870
$ws.Range("J122").Value = 2748.75
$ws.Range("L122").Value = 8246.25
$ws.Range("N122").Value = -13146.25

$ws.Range("H136").Value = 1256056.9
$ws.Range("I136").Value = 4340.073
$ws.Range("J136").Value = 9809455
$ws.Range("K136").Value = 13020.219
$ws.Range("L136").Value = 29428365
$ws.Range("M136").Value = -10470.219
$ws.Range("N136").Value = -29433465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 1212
$ws.Range("J25").Value = 2500
$ws.Range("L25").Value = 2500
$ws.Range("N25").Value = -2970

$ws.Range("H33").Value = 6149.8335
$ws.Range("I33").Value = 1633
$ws.Range("J33").Value = 10666.667
$ws.Range("K33").Value = 1633
$ws.Range("L33").Value = 10666.667
$ws.Range("M33").Value = -1297
$ws.Range("N33").Value = -11338.667

$ws.Range("H86").Value = 5696.4287
$ws.Range("I86").Value = 2306.2222
$ws.Range("K86").Value = 2306.2222
$ws.Range("M86").Value = -1183.2222

$ws.Range("H89").Value = 5696.4287
$ws.Range("I89").Value = 2306.2222
$ws.Range("K89").Value = 11531.111
$ws.Range("M89").Value = -5915.111000000001

$ws.Range("H94").Value = 4353.778
$ws.Range("I94").Value = 3926.4285
$ws.Range("K94").Value = 3926.4285
$ws.Range("M94").Value = -3475.4285

$ws.Range("H134").Value = 21953932
$ws.Range("I134").Value = 1662.6538
$ws.Range("J134").Value = 60004532
$ws.Range("K134").Value = 4987.9614
$ws.Range("L134").Value = 180013596
$ws.Range("M134").Value = -2452.9614
$ws.Range("N134").Value = -180018666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3239.842
$ws.Range("J16").Value = 6915.2
$ws.Range("L16").Value = 6915.2
$ws.Range("N16").Value = -7489.2

$ws.Range("H29").Value = 9999
$ws.Range("J29").Value = 9999
$ws.Range("L29").Value = 9999
$ws.Range("N29").Value = -10585

$ws.Range("H58").Value = 2942.5557
$ws.Range("I58").Value = 2171
$ws.Range("K58").Value = 2171
$ws.Range("M58").Value = -1968

$ws.Range("H62").Value = 3492.818
$ws.Range("J62").Value = 2850
$ws.Range("L62").Value = 2850
$ws.Range("N62").Value = -4098

$ws.Range("H65").Value = 3492.818
$ws.Range("J65").Value = 2850
$ws.Range("L65").Value = 14250
$ws.Range("N65").Value = -20490

$ws.Range("H113").Value = 3239.842
$ws.Range("J113").Value = 6915.2
$ws.Range("L113").Value = 6915.2
$ws.Range("N113").Value = -11255.2

$ws.Range("H132").Value = 12822969
$ws.Range("I132").Value = 2457.125
$ws.Range("J132").Value = 166669120
$ws.Range("K132").Value = 7371.375
$ws.Range("L132").Value = 500007360
$ws.Range("M132").Value = -4841.375
$ws.Range("N132").Value = -500012420

$ws.Range("H134").Value = 1651.7
$ws.Range("I134").Value = 1061.8096
$ws.Range("J134").Value = 3028.111
$ws.Range("K134").Value = 3185.4288
$ws.Range("L134").Value = 9084.332999999999
$ws.Range("M134").Value = -650.4288000000001
$ws.Range("N134").Value = -14154.333

$ws.Range("H136").Value = 2942.5557
$ws.Range("I136").Value = 2171
$ws.Range("K136").Value = 6513
$ws.Range("M136").Value = -3963

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4538830
$ws.Range("I4").Value = 5850104.5
$ws.Range("J4").Value = 167916.17
$ws.Range("K4").Value = 17550313.5
$ws.Range("L4").Value = 503748.51
$ws.Range("M4").Value = -17550201.5
$ws.Range("N4").Value = -503972.51

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = $null

$ws.Range("H122").Value = 6013640
$ws.Range("J122").Value = 4763103.5
$ws.Range("L122").Value = 42867931.5
$ws.Range("N122").Value = -42872831.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 823.8
$ws.Range("I97").Value = 774
$ws.Range("K97").Value = 774
$ws.Range("M97").Value = -278

$ws.Range("H126").Value = 2876.75
$ws.Range("J126").Value = 3171.3333
$ws.Range("L126").Value = 9513.999899999999
$ws.Range("N126").Value = -14453.9999

$ws.Range("H132").Value = 55131210
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 55131210
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 165393630
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -165398690

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8834.666999999999
$ws.Range("I7").Value = 3204.4
$ws.Range("K7").Value = 3204.4
$ws.Range("M7").Value = -3092.4

$ws.Range("H14").Value = 171034.5
$ws.Range("I14").Value = 171034.5
$ws.Range("K14").Value = 171034.5
$ws.Range("M14").Value = -170862.5

$ws.Range("H21").Value = 8000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null

$ws.Range("H43").Value = 2143999.8
$ws.Range("I43").Value = 100000
$ws.Range("K43").Value = 100000
$ws.Range("M43").Value = -99807

$ws.Range("H94").Value = 56249.625

$ws.Range("H126").Value = 8834.666999999999
$ws.Range("I126").Value = 3204.4
$ws.Range("K126").Value = 9613.200000000001
$ws.Range("M126").Value = -7143.200000000001

$ws.Range("H132").Value = 5766.0977
$ws.Range("I132").Value = 2202.8928
$ws.Range("J132").Value = 13440.692
$ws.Range("K132").Value = 6608.678400000001
$ws.Range("L132").Value = 40322.076
$ws.Range("M132").Value = -4078.678400000001
$ws.Range("N132").Value = -45382.076

$ws.Range("H136").Value = 4685.8125
$ws.Range("I136").Value = 4132.778
$ws.Range("J136").Value = 5396.857
$ws.Range("K136").Value = 12398.334
$ws.Range("L136").Value = 16190.571
$ws.Range("M136").Value = -9848.334000000001
$ws.Range("N136").Value = -21290.571

$ws.Range("H141").Value = 159995
$ws.Range("J141").Value = 159995
$ws.Range("L141").Value = 159995
$ws.Range("N141").Value = -170355

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = $null

$ws.Range("H62").Value = 5221.25
$ws.Range("I62").Value = 4963.3335
$ws.Range("J62").Value = 5995
$ws.Range("K62").Value = 4963.3335
$ws.Range("L62").Value = 5995
$ws.Range("M62").Value = -4339.3335
$ws.Range("N62").Value = -7243

$ws.Range("H65").Value = 5221.25
$ws.Range("I65").Value = 4963.3335
$ws.Range("J65").Value = 5995
$ws.Range("K65").Value = 24816.6675
$ws.Range("L65").Value = 29975
$ws.Range("M65").Value = -21696.6675
$ws.Range("N65").Value = -36215

$ws.Range("H132").Value = 32054.303
$ws.Range("I132").Value = 43169.5
$ws.Range("J132").Value = 2413.7778
$ws.Range("K132").Value = 129508.5
$ws.Range("L132").Value = 7241.3334
$ws.Range("M132").Value = -126978.5
$ws.Range("N132").Value = -12301.3334

$ws.Range("H136").Value = 20924.117
$ws.Range("I136").Value = 27885.459
$ws.Range("J136").Value = 2526.2856
$ws.Range("K136").Value = 83656.37699999999
$ws.Range("L136").Value = 7578.8568
$ws.Range("M136").Value = -81106.37699999999
$ws.Range("N136").Value = -12678.8568
